# Apply crypto price/volume updates from the GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold plain-looking numbers (e.g. "233.10") that must stay
# TEXT, matching the source sheet (avoids float round-trip + type drift).
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D17","D19","D20","D23","D26","D27","D28","D29","D32","D34","D38","D39","D43","D44","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.598.48"
$ws.Range("E2").Value = "  +14.10%  "
$ws.Range("D3").Value = "1.825.44"
$ws.Range("E3").Value = "  +8.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "233.10"
$ws.Range("E5").Value = "  +5.21%  "
$ws.Range("E6").Value = "  +5.54%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "32.17"
$ws.Range("E8").Value = "  +7.18%  "
$ws.Range("D9").Value = "46.26"
$ws.Range("E9").Value = "  +4.97%  "
$ws.Range("D10").Value = "0.286"
$ws.Range("E10").Value = "  +8.10%  "
$ws.Range("D11").Value = "0.0683"
$ws.Range("E11").Value = "  +9.68%  "
$ws.Range("D12").Value = "0.0932"
$ws.Range("E12").Value = "  +3.37%  "
$ws.Range("D13").Value = "2.088.43"
$ws.Range("E13").Value = "  +8.61%  "
$ws.Range("D14").Value = "1.825.68"
$ws.Range("E14").Value = "  +8.50%  "
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("D16").Value = "34.529.09"
$ws.Range("E16").Value = "  +13.79%  "
$ws.Range("D17").Value = "10.41"
$ws.Range("E17").Value = "  -2.98%  "
$ws.Range("E18").Value = "  +8.09%  "
$ws.Range("D19").Value = "71.32"
$ws.Range("E19").Value = "  +8.37%  "
$ws.Range("D20").Value = "263.15"
$ws.Range("E20").Value = "  +6.49%  "
$ws.Range("E21").Value = "  +5.73%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "10.58"
$ws.Range("E23").Value = "  +3.77%  "
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "162.12"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").Value = "17.01"
$ws.Range("E27").Value = "  +7.13%  "
$ws.Range("D28").Value = "0.117"
$ws.Range("E28").Value = "  +5.53%  "
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +10.54%  "
$ws.Range("D32").Value = "0.0518"
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("D34").Value = "3.60"
$ws.Range("E34").Value = "  +8.67%  "
$ws.Range("D35").Value = "1.589.82"
$ws.Range("E35").Value = "  +5.97%  "
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("E37").Value = "  +3.52%  "
$ws.Range("D38").Value = "86.39"
$ws.Range("E38").Value = "  +9.74%  "
$ws.Range("D39").Value = "0.635"
$ws.Range("E39").Value = "  +8.36%  "
$ws.Range("E40").Value = "  +4.89%  "
$ws.Range("E41").Value = "  +2.87%  "
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").Value = "0.925"
$ws.Range("E43").Value = "  +8.39%  "
$ws.Range("D44").Value = "2.15"
$ws.Range("E44").Value = "  +7.24%  "
$ws.Range("E45").Value = "  +4.36%  "
$ws.Range("E46").Value = "  +6.41%  "
$ws.Range("D47").Value = "1.977.58"
$ws.Range("E47").Value = "  +8.71%  "
$ws.Range("D48").Value = "54.35"
$ws.Range("E48").Value = "  +4.88%  "
$ws.Range("D49").Value = "5.77"
$ws.Range("E49").Value = "  +6.20%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0125"
$ws.Range("E51").Value = "  +7.13%  "

# Restore default cell style now that the text value is committed, so
# the only observable change is the text content (no style diff).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
